# Insert a new data row at row 564 (pushing the existing rows 564-632 down
# to 565-633), then populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(564).Insert()

$ws.Range("A564").Value2 = 9
$ws.Range("B564").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C564").Value2 = "Metropolitana"
$ws.Range("D564").Value2 = 45212
$ws.Range("E564").Value2 = 13
$ws.Range("F564").Value2 = 100112052
$ws.Range("G564").Value2 = "Albahaca"
$ws.Range("H564").Value2 = "Sin especificar"
$ws.Range("I564").Value2 = "Primera"
$ws.Range("J564").Value2 = 160
$ws.Range("K564").Value2 = 5000
$ws.Range("L564").Value2 = 5000
$ws.Range("M564").Value2 = 5000
$ws.Range("N564").Value2 = "`$/docena de matas"
$ws.Range("O564").Value2 = "Provincia de Chacabuco"
$ws.Range("P564").Value2 = 833
$ws.Range("Q564").Value2 = 6
$ws.Range("R564").Value2 = "Hortaliza"
